$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10504.188787249
$ws.Range("F2").Value = 4.41686729536362

$ws.Range("C3").Value = 9731.12115582882
$ws.Range("F3").Value = 230.941445691466

$ws.Range("C4").Value = 6650.29614570965
$ws.Range("F4").Value = 86.0736358988374

$ws.Range("C5").Value = 6859.70618057313
$ws.Range("F5").Value = 98.7361684428167

$ws.Range("C6").Value = 10574.8991279513
$ws.Range("F6").Value = 277.397400044322

$ws.Range("C7").Value = 9692.22992303822
$ws.Range("F7").Value = 271.494431806375
